# Scientific Integrity Audit: Fixed gene count 617->616
# Replace every standalone occurrence of the (incorrect) gene count "617"
# with the corrected count "616" throughout the manuscript body
# (abstract, intro, results heading/text, figure captions, discussion).

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

# MatchWholeWord ensures only the standalone number "617" is matched
# (not as part of a longer digit run), and ReplaceAll sweeps every
# occurrence across the whole document body in one pass.
$d.Content.Find.Execute(
    "617",      # FindText
    $true,      # MatchCase
    $true,      # MatchWholeWord
    $false,     # MatchWildcards
    $false,     # MatchSoundsLike
    $false,     # MatchAllWordForms
    $true,      # Forward
    1,          # Wrap (wdFindContinue)
    $false,     # Format
    "616",      # ReplaceWith
    2           # Replace (wdReplaceAll)
)
